$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = 'Other found locations'
$ws.Range("I2").Value = ''
$ws.Range("E3").Value = '[Zhongliang%Wang%NULL%0,  Bohan%Yang%NULL%1,  Qianwen%Li%NULL%1,  Lu%Wen%NULL%1,  Ruiguang%Zhang%zrg27@163.com%1]'
$ws.Range("I3").Value = '_PMC'
$ws.Range("E4").Value = '[Ke%Wang%NULL%2,  Zhiguo%Zhang%NULL%2,  Zhiguo%Zhang%NULL%0,  Muqing%Yu%NULL%1,  Yu%Tao%NULL%1,  Min%Xie%xie_m@126.com%1]'
$ws.Range("I4").Value = '_PMC_Springer'
$ws.Range("E5").Value = '[Jianfeng%Wu%NULL%1,  Jianqiang%Huang%NULL%1,  Guochao%Zhu%NULL%1,  Yihao%Liu%NULL%1,  Han%Xiao%NULL%1,  Qian%Zhou%NULL%1,  Xiang%Si%NULL%1,  Hui%Yi%NULL%1,  Cuiping%Wang%NULL%1,  Daya%Yang%NULL%1,  Shuling%Chen%NULL%1,  Xin%Liu%NULL%1,  Zelong%Liu%NULL%1,  Qiongya%Wang%NULL%1,  Qingquan%Lv%NULL%1,  Ying%Huang%NULL%1,  Yang%Yu%NULL%1,  Xiangdong%Guan%NULL%1,  Yanbing%Li%NULL%1,  Krishnarajah%Nirantharakumar%NULL%2,  KarKeung%Cheng%NULL%1,  Sui%Peng%pengsui@vip.163.com%1,  Haipeng%Xiao%xiaohp@mail.sysu.edu.cn%1]'
$ws.Range("I5").Value = '_PMC'
$ws.Range("E6").Value = '[Kaijin%Xu%NULL%1,  Yanfei%Chen%NULL%1,  Jing%Yuan%NULL%2,  Ping%Yi%NULL%1,  Cheng%Ding%NULL%1,  Wenrui%Wu%NULL%2,  Yongtao%Li%NULL%1,  Qin%Ni%NULL%1,  Rongrong%Zou%NULL%1,  Xiaohe%Li%NULL%1,  Min%Xu%NULL%1,  Ying%Zhang%NULL%1,  Hong%Zhao%NULL%1,  Xuan%Zhang%NULL%1,  Liang%Yu%NULL%1,  Junwei%Su%NULL%1,  Guanjing%Lang%NULL%1,  Jun%Liu%NULL%2,  Xiaoxin%Wu%NULL%1,  Yongzheng%Guo%NULL%1,  Jingjing%Tao%NULL%1,  Ding%Shi%NULL%1,  Ling%Yu%NULL%1,  Qing%Cao%NULL%1,  Bing%Ruan%NULL%1,  Lei%Liu%NULL%4,  Zhaoqin%Wang%NULL%2,  Yan%Xu%NULL%1,  Yingxia%Liu%NULL%0,  Jifang%Sheng%NULL%3,  Lanjuan%Li%ljli@zju.edu.cn%4]'
$ws.Range("I6").Value = '_PMC'
$ws.Range("E7").Value = '[Stephen Su%Yang%stephen.yang@mail.mcgill.ca%1,  Jed%Lipes%NULL%2,  Jed%Lipes%NULL%0]'
$ws.Range("I7").Value = '_PMC_Springer'
$ws.Range("E8").Value = '[Lei%Zha%NULL%2,  Shirong%Li%NULL%1,  Lingling%Pan%NULL%1,  Boris%Tefsen%NULL%1,  Yeshan%Li%NULL%1,  Neil%French%NULL%1,  Liyun%Chen%NULL%1,  Gang%Yang%NULL%1,  Elmer V%Villanueva%Villanueva@xjtlu.edu.cn%1]'
$ws.Range("I8").Value = '_PMC'
$ws.Range("E9").Value = '[Francesco%Salton%NULL%1,  Paola%Confalonieri%NULL%1,  G Umberto%Meduri%NULL%1,  Pierachille%Santus%NULL%1,  Sergio%Harari%NULL%1,  Raffaele%Scala%NULL%1,  Simone%Lanini%NULL%1,  Valentina%Vertui%NULL%1,  Tiberio%Oggionni%NULL%1,  Antonella%Caminati%NULL%1,  Vincenzo%Patruno%NULL%1,  Mario%Tamburrini%NULL%1,  Alessandro%Scartabellati%NULL%1,  Mara%Parati%NULL%1,  Massimiliano%Villani%NULL%1,  Dejan%Radovanovic%NULL%1,  Sara%Tomassetti%NULL%1,  Claudia%Ravaglia%NULL%1,  Venerino%Poletti%NULL%1,  Andrea%Vianello%NULL%1,  Anna Talia%Gaccione%NULL%1,  Luca%Guidelli%NULL%1,  Rita%Raccanelli%NULL%1,  Paolo%Lucernoni%NULL%1,  Donato%Lacedonia%NULL%1,  Maria Pia%Foschino Barbaro%NULL%1,  Stefano%Centanni%NULL%1,  Michele%Mondoni%NULL%1,  Matteo%Davì%NULL%1,  Alberto%Fantin%NULL%1,  Xueyuan%Cao%NULL%1,  Lucio%Torelli%NULL%1,  Antonella%Zucchetto%NULL%1,  Marcella%Montico%NULL%1,  Annalisa%Casarin%NULL%1,  Micaela%Romagnoli%NULL%1,  Stefano%Gasparini%NULL%1,  Martina%Bonifazi%NULL%1,  Pierlanfranco%D’Agaro%NULL%1,  Alessandro%Marcello%NULL%1,  Danilo%Licastro%NULL%1,  Barbara%Ruaro%NULL%1,  Maria Concetta%Volpe%NULL%1,  Reba%Umberger%NULL%1,  Marco%Confalonieri%mconfalonieri@units.it%1]'
$ws.Range("I9").Value = '_PMC'
$ws.Range("E10").Value = '[Yinzhong%Shen%NULL%2,  Fang%Zheng%NULL%2,  Fang%Zheng%NULL%0,  Danfeng%Sun%NULL%1,  Yun%Ling%NULL%5,  Jun%Chen%NULL%2,  Feng%Li%NULL%3,  Feng%Li%NULL%0,  Tao%Li%NULL%3,  Zhiping%Qian%NULL%2,  Yuyi%Zhang%NULL%2,  Qingnian%Xu%NULL%2,  Li%Liu%NULL%2,  Qin%Huang%NULL%1,  Fei%Shan%NULL%2,  Lie%Xu%NULL%1,  Jun%Wu%NULL%1,  Zhaoqin%Zhu%NULL%1,  Zhigang%Song%NULL%2,  Shenyang%Li%NULL%1,  Yuxin%Shi%NULL%4,  Jianliang%Zhang%NULL%1,  Xueyun%Wu%NULL%1,  Joshua B.%Mendelsohn%NULL%2,  Joshua B.%Mendelsohn%NULL%0,  Tongyu%Zhu%NULL%2,  Hongzhou%Lu%NULL%5]'
$ws.Range("I10").Value = '_PMC'
$ws.Range("E11").Value = '[José Luis%Callejas Rubio%jlcalleja@telefonica.net%1,  Juan de Dios%Luna del Castillo%NULL%1,  Javier%de la Hera Fernández%NULL%1,  Emilio%Guirao Arrabal%NULL%1,  Manuel%Colmenero Ruiz%NULL%1,  Norberto%Ortego Centeno%NULL%1]'
$ws.Range("I11").Value = '_PMC_elsevier'
$ws.Range("E12").Value = '[Ding%Shi%NULL%1,  Wenrui%Wu%NULL%0,  Qing%Wang%NULL%1,  Kaijin%Xu%NULL%1,  Jiaojiao%Xie%NULL%1,  Jingjing%Wu%NULL%1,  Longxian%Lv%NULL%1,  Jifang%Sheng%NULL%0,  Jing%Guo%NULL%1,  Kaicen%Wang%NULL%1,  Daiqiong%Fang%NULL%1,  Yating%Li%NULL%1,  Lanjuan%Li%ljli@zju.edu.cn%0]'
$ws.Range("I12").Value = '_PMC'
$ws.Range("E13").Value = '[Yin%Wang%NULL%2,  Weiwei%Jiang%NULL%1,  Qi%He%NULL%1,  Cheng%Wang%NULL%1,  Baoju%Wang%NULL%1,  Pan%Zhou%NULL%1,  Nianguo%Dong%dongnianguo63@gmail.com%1,  Qiaoxia%Tong%2013xh0859@hust.edu.cn%1]'
$ws.Range("I13").Value = '_PMC_Springer'
$ws.Range("E14").Value = '[Yang%Li%NULL%3,  Xian%Zhou%NULL%1,  Tao%Li%NULL%0,  Shiji%Chan%NULL%1,  Yiqi%Yu%NULL%1,  Jing-Wen%Ai%NULL%1,  Haocheng%Zhang%NULL%1,  Feng%Sun%NULL%1,  Qiran%Zhang%NULL%1,  Lei%Zhu%NULL%1,  Lingyun%Shao%NULL%1,  Bin%Xu%NULL%1,  Wenhong%Zhang%NULL%1]'
$ws.Range("I14").Value = '_PMC'
$ws.Range("E15").Value = '[Jing%Liu%NULL%1,  Xiaobin%Zheng%NULL%1,  Yiying%Huang%NULL%1,  Hong%Shan%shanhong@mail.sysu.edu.cn%5,  Jin%Huang%hjin@mail.sysu.edu.cn%1]'
$ws.Range("I15").Value = '_PMC_elsevier'
$ws.Range("I16").Value = ''
$ws.Range("E17").Value = '[Xiang-Hong%Yang%NULL%1,  Ran-Ran%Li%NULL%1,  Ren-Hua%Sun%NULL%1,  Jiao%Liu%NULL%1,  De-Chang%Chen%NULL%1,  Peng%Lyu%NULL%2,  Peng%Lyu%NULL%0]'
$ws.Range("I17").Value = '_PMC'
$ws.Range("E18").Value = '[Xiaofan%Lu%NULL%2,  Taige%Chen%NULL%1,  Yang%Wang%NULL%1,  Jun%Wang%NULL%1,  Fangrong%Yan%f.r.yan@163.com%1]'
$ws.Range("I18").Value = '_PMC_Springer'
$ws.Range("E19").Value = '[Yiming%Ma%NULL%1,  Huihui%Zeng%NULL%1,  Zijie%Zhan%NULL%1,  Huanhuan%Lu%NULL%1,  Zihang%Zeng%NULL%1,  Chenjie%He%NULL%1,  Xiangming%Liu%NULL%1,  Chen%Chen%NULL%1,  Qingwu%Qin%NULL%1,  Jia%He%NULL%1,  Zhiguo%Zhou%NULL%1,  Peng%Huang%NULL%2,  Mingyan%Jiang%NULL%1,  Dingding%Deng%NULL%1,  Xin%Liao%NULL%1,  Zhi%Xiang%NULL%1,  Xiaoying%Huang%NULL%1,  Yan%Chen%NULL%1,  Ping%Chen%NULL%1]'
$ws.Range("I19").Value = '_PMC'
$ws.Range("E20").Value = '[Monil%Majmundar%NULL%1,  Tikal%Kansara%NULL%2,  Tikal%Kansara%NULL%0,  Joanna Marta%Lenik%NULL%1,  Hansang%Park%NULL%1,  Kuldeep%Ghosh%NULL%1,  Rajkumar%Doshi%NULL%1,  Palak%Shah%NULL%1,  Ashish%Kumar%NULL%1,  Hossam%Amin%NULL%2,  Hossam%Amin%NULL%0,  Shobhana%Chaudhari%NULL%1,  Imnett%Habtes%NULL%1,  Muhammad%Adrish%NULL%8,  Muhammad%Adrish%NULL%0]'
$ws.Range("I20").Value = '_PMC'
$ws.Range("E21").Value = '[Malgorzata%Mikulska%NULL%1,  Laura Ambra%Nicolini%NULL%2,  Laura Ambra%Nicolini%NULL%0,  Alessio%Signori%NULL%1,  Antonio%Di Biagio%NULL%1,  Chiara%Sepulcri%NULL%2,  Chiara%Sepulcri%NULL%0,  Chiara%Russo%NULL%1,  Silvia%Dettori%NULL%1,  Marco%Berruti%NULL%1,  Maria Pia%Sormani%NULL%1,  Daniele Roberto%Giacobbe%NULL%1,  Antonio%Vena%NULL%2,  Antonio%Vena%NULL%0,  Andrea%De Maria%NULL%1,  Chiara%Dentone%NULL%1,  Lucia%Taramasso%NULL%1,  Michele%Mirabella%NULL%2,  Michele%Mirabella%NULL%0,  Laura%Magnasco%NULL%1,  Sara%Mora%NULL%1,  Emanuele%Delfino%NULL%1,  Federica%Toscanini%NULL%1,  Elisa%Balletto%NULL%1,  Anna Ida%Alessandrini%NULL%1,  Federico%Baldi%NULL%1,  Federica%Briano%NULL%1,  Marco%Camera%NULL%1,  Ferdinando%Dodi%NULL%1,  Antonio%Ferrazin%NULL%1,  Laura%Labate%NULL%1,  Giovanni%Mazzarello%NULL%1,  Rachele%Pincino%NULL%1,  Federica%Portunato%NULL%1,  Stefania%Tutino%NULL%1,  Emanuela%Barisione%NULL%1,  Bianca%Bruzzone%NULL%1,  Andrea%Orsi%NULL%1,  Eva%Schenone%NULL%1,  Nirmala%Rosseti%NULL%1,  Elisabetta%Sasso%NULL%1,  Giorgio%Da Rin%NULL%1,  Paolo%Pelosi%NULL%1,  Sabrina%Beltramini%NULL%1,  Mauro%Giacomini%NULL%1,  Giancarlo%Icardi%NULL%2,  Giancarlo%Icardi%NULL%0,  Angelo%Gratarola%NULL%1,  Matteo%Bassetti%NULL%1,  Muhammad%Adrish%NULL%0,  Muhammad%Adrish%NULL%0,  Muhammad%Adrish%NULL%0,  Muhammad%Adrish%NULL%0,  Muhammad%Adrish%NULL%0,  Muhammad%Adrish%NULL%0]'
$ws.Range("I21").Value = '_PMC'
$ws.Range("E22").Value = '[Mao%Huang%NULL%2,  Yi%Yang%NULL%1,  Futai%Shang%NULL%1,  Yishan%Zheng%NULL%1,  Wenjing%Zhao%NULL%1,  Liang%Luo%NULL%1,  Xudong%Han%NULL%1,  Aihua%Lin%NULL%1,  Hongsheng%Zhao%NULL%1,  Qing%Gu%NULL%1,  Yi%Shi%NULL%1,  Jun%Li%NULL%1,  Xingxiang%Xu%NULL%1,  Kexi%Liu%NULL%1,  YiJun%Deng%NULL%1,  Quan%Cao%NULL%1,  Weiwei%Wang%NULL%1]'
$ws.Range("I22").Value = '_PMC_elsevier'
$ws.Range("E23").Value = '[Yun%Feng%NULL%4,  Yun%Ling%NULL%0,  Yun%Ling%NULL%0,  Tao%Bai%NULL%2,  Tao%Bai%NULL%0,  Yusang%Xie%NULL%1,  Jie%Huang%NULL%2,  Jie%Huang%NULL%0,  Jian%Li%NULL%1,  Weining%Xiong%NULL%1,  Dexiang%Yang%NULL%1,  Rong%Chen%NULL%2,  Fangying%Lu%NULL%1,  Yunfei%Lu%NULL%1,  Xuhui%Liu%NULL%1,  Yuqing%Chen%NULL%2,  Yuqing%Chen%NULL%0,  Xin%Li%NULL%1,  Yong%Li%NULL%1,  Hanssa Dwarka%Summah%NULL%1,  Huihuang%Lin%NULL%1,  Jiayang%Yan%NULL%1,  Min%Zhou%NULL%5,  Hongzhou%Lu%NULL%0,  Hongzhou%Lu%NULL%0,  Jieming%Qu%NULL%5,  Jieming%Qu%NULL%0]'
$ws.Range("I23").Value = '_PMC'
$ws.Range("E24").Value = '[Ana%Fernández-Cruz%NULL%1,  Belén%Ruiz-Antorán%NULL%2,  Belén%Ruiz-Antorán%NULL%0,  Ana%Muñoz-Gómez%NULL%1,  Aránzazu%Sancho-López%NULL%1,  Patricia%Mills-Sánchez%NULL%1,  Gustavo Adolfo%Centeno-Soto%NULL%1,  Silvia%Blanco-Alonso%NULL%1,  Laura%Javaloyes-Garachana%NULL%1,  Amy%Galán-Gómez%NULL%1,  Ángela%Valencia-Alijo%NULL%1,  Javier%Gómez-Irusta%NULL%1,  Concepción%Payares-Herrera%NULL%1,  Ignacio%Morrás-Torre%NULL%1,  Enrique%Sánchez-Chica%NULL%1,  Laura%Delgado-Téllez-de-Cepeda%NULL%1,  Alejandro%Callejas-Díaz%NULL%1,  Antonio%Ramos-Martínez%NULL%1,  Elena%Múñez-Rubio%NULL%2,  Elena%Múñez-Rubio%NULL%0,  Cristina%Avendaño-Solá%NULL%1]'
$ws.Range("I24").Value = '_PMC'
$ws.Range("E25").Value = '[Carmine%Gazzaruso%c.gazzaruso@gmail.com%1,  Nicoletta%Carlo Stella%NULL%2,  Nicoletta%Carlo Stella%NULL%0,  Giuseppe%Mariani%NULL%1,  Anna%Tamburlini%NULL%1,  Pietro%Garini%NULL%1,  Elena%Freddi%NULL%1,  Carolina%Ravetto%NULL%1,  Adriana%Coppola%NULL%1,  Pietro%Gallotti%NULL%1]'
$ws.Range("I25").Value = '_PMC_Springer'
$ws.Range("E26").Value = '[Yan%Hu%NULL%2,  Tingping%Wang%NULL%2,  Zhimin%Hu%NULL%2,  Xuan%Wang%NULL%2,  Zhengbin%Zhang%NULL%2,  Li%Li%NULL%9,  Peng%Peng%NULL%5]'
$ws.Range("I26").Value = '_PMC_elsevier'
$ws.Range("I27").Value = '_PMC'
$ws.Range("E28").Value = '[Yan%Hu%NULL%0,  Tingping%Wang%NULL%0,  Zhimin%Hu%NULL%0,  Xuan%Wang%NULL%0,  Zhengbin%Zhang%NULL%0,  Li%Li%NULL%0,  Peng%Peng%NULL%0]'
$ws.Range("I28").Value = '_PMC_elsevier'
$ws.Range("E29").Value = '[Jianlei%Cao%NULL%0,  Wen-Jun%Tu%tuwenjun@irm-cams.ac.cn%2,  Wenlin%Cheng%NULL%2,  Lei%Yu%NULL%2,  Ya-Kun%Liu%NULL%2,  Xiaoyong%Hu%NULL%2,  Qiang%Liu%NULL%2]'
$ws.Range("I29").Value = '_PMC'
$ws.Range("I30").Value = ''
$ws.Range("E31").Value = '[Firouzé%Bani-Sadr%NULL%1,  Maxime%Hentzien%NULL%1,  Madeline%Pascard%NULL%1,  Yohan%N''Guyen%NULL%1,  Amélie%Servettaz%NULL%1,  Laurent%Andreoletti%NULL%1,  Lukshe%Kanagaratnam%NULL%1,  Damien%Jolly%NULL%1]'
$ws.Range("I31").Value = '_PMC_elsevier'
$ws.Range("I32").Value = ''
$ws.Range("E33").Value = '[Raef%Fadel%NULL%2,  Austin R%Morrison%NULL%1,  Amit%Vahia%NULL%1,  Zachary R%Smith%NULL%1,  Zohra%Chaudhry%NULL%1,  Pallavi%Bhargava%NULL%1,  Joseph%Miller%NULL%1,  Rachel M%Kenney%NULL%1,  George%Alangaden%NULL%1,  Mayur S%Ramesh%mramesh1@hfhs.org%1,  Varidhi%Nauriyal%NULL%1,  Jayanth%Lakshmikanth%NULL%1,  Asif%Abdul Hamed%NULL%1,  Owais%Nadeem%NULL%1,  Kristin%Griebe%NULL%1,  Joseph M%Johnson%NULL%1,  Patrick%Bradley%NULL%1,  Junior%Uduman%NULL%1,  Sara%Hegab%NULL%1,  Jennifer%Swiderek%NULL%1,  Amanda%Godfrey%NULL%1,  Jeffrey%Jennings%NULL%1,  Jayna%Gardner-Gray%NULL%1,  Adam%Ackerman%NULL%1,  Jonathan%Lezotte%NULL%1,  Joseph%Ruhala%NULL%1,  Linoj%Samuel%NULL%1,  Robert J%Tibbetts%NULL%1,  Indira%Brar%NULL%1,  John%McKinnon%NULL%1,  Geehan%Suleyman%NULL%1,  Nicholas%Yared%NULL%1,  Erica%Herc%NULL%1,  Jonathan%Williams%NULL%1,  Odaliz Abreu%Lanfranco%NULL%1,  Anne%Chen%NULL%1,  Marcus%Zervos%NULL%1,  Eric%Scher%NULL%1]'
$ws.Range("I33").Value = '_PMC'
$ws.Range("E34").Value = '[Tomasz%Chroboczek%tchroboczek@gmail.com%1,  Marie%Lacoste%NULL%2,  Marie%Lacoste%NULL%0,  Chloe%Wackenheim%NULL%1,  Thibaut%Challan-Belval%NULL%1,  Benjamin%Amar%NULL%1,  Thomas%Boisson%NULL%1,  Jason%Hubac%NULL%1,  Dominique%Leduc%NULL%1,  Colleen%Masse%NULL%1,  Victor%Dechaene%NULL%1,  Laetitia%Touihri-Maximin%NULL%1,  Sandrine%Megessier%NULL%1,  Camille%Lassale%NULL%1]'
$ws.Range("I34").Value = '_PMC'
$ws.Range("I35").Value = ''
$ws.Range("E36").Value = '[Christiane Maria Prado%Jeronimo%NULL%2,  Maria Eduarda Leão%Farias%NULL%1,  Fernando Fonseca Almeida%Val%NULL%1,  Vanderson Souza%Sampaio%NULL%1,  Marcia Almeida Araújo%Alexandre%NULL%1,  Gisely Cardoso%Melo%NULL%1,  Izabella Picinin%Safe%NULL%1,  Mayla Gabriela Silva%Borba%NULL%1,  Rebeca Linhares%Abreu-Netto%NULL%1,  Alex Bezerra Silva%Maciel%NULL%1,  João Ricardo Silva%Neto%NULL%1,  Lucas Barbosa%Oliveira%NULL%1,  Erick Frota Gomes%Figueiredo%NULL%1,  Kelry Mazurega Oliveira%Dinelly%NULL%1,  Maria Gabriela de Almeida%Rodrigues%NULL%1,  Marcelo%Brito%NULL%1,  Maria Paula Gomes%Mourão%NULL%1,  Guilherme Augusto%Pivoto João%NULL%1,  Ludhmila Abrahão%Hajjar%NULL%1,  Quique%Bassat%NULL%1,  Gustavo Adolfo Sierra%Romero%NULL%1,  Felipe Gomes%Naveca%NULL%1,  Heline Lira%Vasconcelos%NULL%1,  Michel de Araújo%Tavares%NULL%1,  José Diego%Brito-Sousa%NULL%1,  Fabio Trindade Maranhão%Costa%NULL%1,  Maurício Lacerda%Nogueira%NULL%1,  Djane%Baía-da-Silva%NULL%1,  Mariana Simão%Xavier%NULL%1,  Wuelton Marcelo%Monteiro%NULL%1,  Marcus Vinícius Guimarães%Lacerda%marcuslacerda.br@gmail.com%1,  NULL%NULL%NULL%22]'
$ws.Range("I36").Value = '_PMC'
$ws.Range("E37").Value = '[Qiang%Li%NULL%1,  Weixia%Li%NULL%1,  Yinpeng%Jin%NULL%1,  Wei%Xu%NULL%1,  Chenlu%Huang%NULL%1,  Li%Li%NULL%0,  Yuxian%Huang%NULL%1,  Qingchun%Fu%fuqingchun@shphc.org.cn%1,  Liang%Chen%chenliang@shphc.org.cn%1]'
$ws.Range("I37").Value = '_PMC_Springer'
$ws.Range("E38").Value = '[Sijia%Li%NULL%2,  Zhigang%Hu%NULL%1,  Xinyu%Song%songxinyu@ctgu.edu.cn%2,  Xinyu%Song%songxinyu@ctgu.edu.cn%0]'
$ws.Range("I38").Value = '_PMC'
$ws.Range("E39").Value = '[Brian C%Nelson%brn9016@nyp.org%1,  Justin%Laracy%NULL%1,  Sherif%Shoucri%NULL%1,  Donald%Dietz%NULL%1,  Jason%Zucker%NULL%1,  Nina%Patel%NULL%1,  Magdalena E%Sobieszczyk%NULL%2,  Christine J%Kubin%NULL%1,  Angela%Gomez-Simmonds%NULL%1]'
$ws.Range("I39").Value = '_PMC_elsevier'
$ws.Range("I40").Value = ''
